$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - copy formatting from H1 (same header style), then set values
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J8
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 9

$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 9

$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 6
